# Weekly update: insert 4 new price rows (newest week) above the existing
# data block in the "Piña" sheet, pushing all prior rows down by 4.
# The repeating 4-row block pattern is:
#   Calidad: Especial / Primera / Segunda / Tercera
#   Unidad:  $/caja 10 / 12 / 14 / 16 unidades
#   Kg/unidad (T): 10 / 12 / 14 / 16
# with fixed constant columns A, B, C, E, F, G, H, I, J, K, R for every row
# in this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows right before the current row 160, shifting the rest
# of the table (old rows 160:199) down to 164:203.
$ws.Rows("160:163").Insert()

# --- Row 160 (Especial) ---
$ws.Range("A160").Value = 1
$ws.Range("B160").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C160").Value = "Arica y Parinacota"
$ws.Range("D160").Value = 44722
$ws.Range("E160").Value = 15
$ws.Range("F160").Value = "Fruta"
$ws.Range("G160").Value = 100108
$ws.Range("H160").Value = "Tropicales y subtropicales"
$ws.Range("I160").Value = 100108005
$ws.Range("J160").Value = "Piña"
$ws.Range("K160").Value = "Caramelo"
$ws.Range("L160").Value = "Especial"
$ws.Range("M160").Value = 200
$ws.Range("N160").Value = 20000
$ws.Range("O160").Value = 21000
$ws.Range("P160").Value = 20500
$ws.Range("Q160").Value = "$/caja 10 unidades"
$ws.Range("R160").Value = "Ecuador"
$ws.Range("S160").Value = 2050
$ws.Range("T160").Value = 10

# --- Row 161 (Primera) ---
$ws.Range("A161").Value = 1
$ws.Range("B161").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C161").Value = "Arica y Parinacota"
$ws.Range("D161").Value = 44722
$ws.Range("E161").Value = 15
$ws.Range("F161").Value = "Fruta"
$ws.Range("G161").Value = 100108
$ws.Range("H161").Value = "Tropicales y subtropicales"
$ws.Range("I161").Value = 100108005
$ws.Range("J161").Value = "Piña"
$ws.Range("K161").Value = "Caramelo"
$ws.Range("L161").Value = "Primera"
$ws.Range("M161").Value = 270
$ws.Range("N161").Value = 20000
$ws.Range("O161").Value = 21000
$ws.Range("P161").Value = 20500
$ws.Range("Q161").Value = "$/caja 12 unidades"
$ws.Range("R161").Value = "Ecuador"
$ws.Range("S161").Value = 1708
$ws.Range("T161").Value = 12

# --- Row 162 (Segunda) ---
$ws.Range("A162").Value = 1
$ws.Range("B162").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C162").Value = "Arica y Parinacota"
$ws.Range("D162").Value = 44722
$ws.Range("E162").Value = 15
$ws.Range("F162").Value = "Fruta"
$ws.Range("G162").Value = 100108
$ws.Range("H162").Value = "Tropicales y subtropicales"
$ws.Range("I162").Value = 100108005
$ws.Range("J162").Value = "Piña"
$ws.Range("K162").Value = "Caramelo"
$ws.Range("L162").Value = "Segunda"
$ws.Range("M162").Value = 270
$ws.Range("N162").Value = 20000
$ws.Range("O162").Value = 21000
$ws.Range("P162").Value = 20500
$ws.Range("Q162").Value = "$/caja 14 unidades"
$ws.Range("R162").Value = "Ecuador"
$ws.Range("S162").Value = 1464
$ws.Range("T162").Value = 14

# --- Row 163 (Tercera) ---
$ws.Range("A163").Value = 1
$ws.Range("B163").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C163").Value = "Arica y Parinacota"
$ws.Range("D163").Value = 44722
$ws.Range("E163").Value = 15
$ws.Range("F163").Value = "Fruta"
$ws.Range("G163").Value = 100108
$ws.Range("H163").Value = "Tropicales y subtropicales"
$ws.Range("I163").Value = 100108005
$ws.Range("J163").Value = "Piña"
$ws.Range("K163").Value = "Caramelo"
$ws.Range("L163").Value = "Tercera"
$ws.Range("M163").Value = 270
$ws.Range("N163").Value = 20000
$ws.Range("O163").Value = 21000
$ws.Range("P163").Value = 20500
$ws.Range("Q163").Value = "$/caja 16 unidades"
$ws.Range("R163").Value = "Ecuador"
$ws.Range("S163").Value = 1281
$ws.Range("T163").Value = 16
